$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 7545.303
$ws.Range("I51").Value = 16624.25
$ws.Range("K51").Value = 16624.25
$ws.Range("M51").Value = -16140.25

$ws.Range("H98").Value = 1539.1428
$ws.Range("I98").Value = 1539.1428
$ws.Range("K98").Value = 1539.1428
$ws.Range("M98").Value = -41.14280000000008

$ws.Range("H122").Value = 1539.1428
$ws.Range("I122").Value = 1539.1428
$ws.Range("K122").Value = 4617.428400000001
$ws.Range("M122").Value = -2167.428400000001

$ws.Range("H137").Value = 1669.491
$ws.Range("I137").Value = 1388.2927
$ws.Range("J137").Value = 2493
$ws.Range("K137").Value = 4164.8781
$ws.Range("L137").Value = 7479
$ws.Range("M137").Value = -1614.8781
$ws.Range("N137").Value = -12579

$ws.Range("H138").Value = 4217.905
$ws.Range("J138").Value = 8209.875
$ws.Range("L138").Value = 24629.625
$ws.Range("N138").Value = -34909.625

$ws.Range("H141").Value = 8774265
$ws.Range("I141").Value = 10001843
$ws.Range("K141").Value = 30005529
$ws.Range("M141").Value = -30000349

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 11055872
$ws.Range("I61").Value = 11114532
$ws.Range("J61").Value = 10000000
$ws.Range("K61").Value = 11114532
$ws.Range("L61").Value = 10000000
$ws.Range("M61").Value = -11114320
$ws.Range("N61").Value = -10000424

$ws.Range("H125").Value = 77619.664
$ws.Range("J125").Value = 77619.664
$ws.Range("L125").Value = 77619.664
$ws.Range("N125").Value = -87459.664

$ws.Range("H136").Value = 11055872
$ws.Range("I136").Value = 11114532
$ws.Range("J136").Value = 10000000
$ws.Range("K136").Value = 33343596
$ws.Range("L136").Value = 30000000
$ws.Range("M136").Value = -33341046
$ws.Range("N136").Value = -30005100

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 4763827
$ws.Range("I134").Value = 1861.3158
$ws.Range("K134").Value = 5583.9474
$ws.Range("M134").Value = -3048.9474

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 3595.6
$ws.Range("I58").Value = 3492
$ws.Range("J58").Value = 3714
$ws.Range("K58").Value = 3492
$ws.Range("L58").Value = 3714
$ws.Range("M58").Value = -3289
$ws.Range("N58").Value = -4120

$ws.Range("H132").Value = 2481.1428
$ws.Range("I132").Value = 2478
$ws.Range("K132").Value = 7434
$ws.Range("M132").Value = -4904

$ws.Range("H136").Value = 3595.6
$ws.Range("I136").Value = 3492
$ws.Range("J136").Value = 3714
$ws.Range("K136").Value = 10476
$ws.Range("L136").Value = 11142
$ws.Range("M136").Value = -7926
$ws.Range("N136").Value = -16242

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 90.40000000000001
$ws.Range("I2").Value = 102.57143
$ws.Range("J2").Value = 62
$ws.Range("K2").Value = 615.42858
$ws.Range("L2").Value = 372
$ws.Range("M2").Value = -502.42858
$ws.Range("N2").Value = -598

$ws.Range("H23").Value = 1031.7059
$ws.Range("J23").Value = 1243.6154
$ws.Range("L23").Value = 3730.8462
$ws.Range("N23").Value = -4200.8462

$ws.Range("H37").Value = 92498.8
$ws.Range("J37").Value = 92498.8
$ws.Range("L37").Value = 277496.4
$ws.Range("N37").Value = -277720.4

$ws.Range("H80").Value = 33335094
$ws.Range("I80").Value = 55556892
$ws.Range("J80").Value = 2400
$ws.Range("K80").Value = 166670676
$ws.Range("L80").Value = 7200
$ws.Range("M80").Value = -166669740
$ws.Range("N80").Value = -9072

$ws.Range("H83").Value = 33335094
$ws.Range("I83").Value = 55556892
$ws.Range("J83").Value = 2400
$ws.Range("K83").Value = 500012028
$ws.Range("L83").Value = 21600
$ws.Range("M83").Value = -500007348
$ws.Range("N83").Value = -30960

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2461.8262
$ws.Range("I102").Value = 2213.5881
$ws.Range("K102").Value = 2213.5881
$ws.Range("M102").Value = -591.5880999999999

$ws.Range("H122").Value = 1837208.5
$ws.Range("I122").Value = 2360576.8
$ws.Range("J122").Value = 5420
$ws.Range("K122").Value = 7081730.399999999
$ws.Range("L122").Value = 16260
$ws.Range("M122").Value = -7079280.399999999
$ws.Range("N122").Value = -21160

$ws.Range("H126").Value = 4299.7334
$ws.Range("I126").Value = 3553.1428
$ws.Range("J126").Value = 4953
$ws.Range("K126").Value = 10659.4284
$ws.Range("L126").Value = 14859
$ws.Range("M126").Value = -8189.428400000001
$ws.Range("N126").Value = -19799

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2319.05
$ws.Range("I16").Value = 1540.7222
$ws.Range("J16").Value = 9324
$ws.Range("K16").Value = 1540.7222
$ws.Range("L16").Value = 9324
$ws.Range("M16").Value = -1370.7222
$ws.Range("N16").Value = -9664

$ws.Range("H22").Value = 12008.182
$ws.Range("J22").Value = 2666.3333
$ws.Range("L22").Value = 2666.3333
$ws.Range("N22").Value = -3256.3333

$ws.Range("H27").Value = 12008.182
$ws.Range("J27").Value = 2666.3333
$ws.Range("L27").Value = 2666.3333
$ws.Range("N27").Value = -2880.3333

$ws.Range("H82").Value = 5708.3335
$ws.Range("I82").Value = 3028.4285
$ws.Range("J82").Value = 9460.200000000001
$ws.Range("K82").Value = 3028.4285
$ws.Range("L82").Value = 9460.200000000001
$ws.Range("M82").Value = -2667.4285
$ws.Range("N82").Value = -10182.2

$ws.Range("H85").Value = 5708.3335
$ws.Range("I85").Value = 3028.4285
$ws.Range("J85").Value = 9460.200000000001
$ws.Range("K85").Value = 3028.4285
$ws.Range("L85").Value = 9460.200000000001
$ws.Range("M85").Value = -1780.4285
$ws.Range("N85").Value = -11956.2

$ws.Range("H100").Value = 15627739
$ws.Range("I100").Value = 1701.4445
$ws.Range("J100").Value = 35718360
$ws.Range("K100").Value = 1701.4445
$ws.Range("L100").Value = 35718360
$ws.Range("M100").Value = -1160.4445
$ws.Range("N100").Value = -35719442

$ws.Range("H122").Value = 3939.675
$ws.Range("I122").Value = 3502.4722
$ws.Range("K122").Value = 10507.4166
$ws.Range("M122").Value = -8057.4166

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H106").Value = 98892
$ws.Range("J106").Value = 98892
$ws.Range("L106").Value = 98892
$ws.Range("N106").Value = -101416

$ws.Range("H132").Value = 770828.9399999999
$ws.Range("I132").Value = 1508.875
$ws.Range("K132").Value = 4526.625
$ws.Range("M132").Value = -1996.625

$ws.Range("H136").Value = 233974.98
$ws.Range("I136").Value = 943.5
$ws.Range("J136").Value = 590376.0600000001
$ws.Range("K136").Value = 2830.5
$ws.Range("L136").Value = 1771128.18
$ws.Range("M136").Value = -280.5
$ws.Range("N136").Value = -1776228.18
